# Update OCR word-box table: split the previous 4 multi-word rows into
# one row per individual word (rows 2-15), with refreshed x/y/width/height
# values for each word's bounding box.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ R = 2;  A = "19th ";         B = 486; C = 644;    D = 43;  E = 23 },
    @{ R = 3;  A = "century. ";     B = 529; C = 644;    D = 73;  E = 23 },
    @{ R = 4;  A = "“Romantic ";    B = 752; C = 644;    D = 93;  E = 23 },
    @{ R = 5;  A = "period” ";      B = 845; C = 644;    D = 66;  E = 23 },
    @{ R = 6;  A = "the ";          B = 656; C = 669;    D = 32;  E = 23 },
    @{ R = 7;  A = "Classical ";    B = 688; C = 669;    D = 83;  E = 23 },
    @{ R = 8;  A = "period, ";      B = 771; C = 669;    D = 64;  E = 23 },
    @{ R = 9;  A = "18th-century "; B = 292; C = 735.4;  D = 112; E = 23 },
    @{ R = 10; A = "Pastoral,” ";   B = 455; C = 1034.6; D = 88;  E = 23 },
    @{ R = 11; A = "sea ";          B = 272; C = 1126;   D = 37;  E = 23 },
    @{ R = 12; A = "coming ";       B = 309; C = 1126;   D = 68;  E = 23 },
    @{ R = 13; A = "into ";         B = 377; C = 1126;   D = 37;  E = 23 },
    @{ R = 14; A = "Fingal’s ";     B = 414; C = 1126;   D = 72;  E = 23 },
    @{ R = 15; A = "Cave ";         B = 486; C = 1126;   D = 51;  E = 23 }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
}
